# Update "Update countries & provincias Spain" data
# - Refresh the "last updated" timestamp in A1
# - Update Covid-19 daily numbers for a number of countries
# - Tunez overtakes Uruguay / Bosnia y Herzegovina (re-sorted by total cases)
# - Republica de Yibuti overtakes Tanzania (re-sorted by total cases)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 23:29"

# --- Straightforward numeric updates (country stays in the same row) --
# Row 4: Estados Unidos
$ws.Range("B4").Value = 121043
$ws.Range("C4").Value = 16917
$ws.Range("D4").Value = 3231
$ws.Range("E4").Value = 115792
$ws.Range("G4").Value = 324
$ws.Range("H4").Value = 2020

# Row 6: China
$ws.Range("C6").Value = 0
$ws.Range("G6").Value = 0

# Row 7: España
$ws.Range("B7").Value = 73235
$ws.Range("C7").Value = 7516
$ws.Range("E7").Value = 54968
$ws.Range("G7").Value = 844
$ws.Range("H7").Value = 5982

# Row 13: Paises Bajos
$ws.Range("F13").Value = 914

# Row 16: Austria
$ws.Range("B16").Value = 8271
$ws.Range("C16").Value = 574
$ws.Range("E16").Value = 7978
$ws.Range("F16").Value = 135

# Row 17: Turquia
$ws.Range("F17").Value = 445

# Row 18: Canada
$ws.Range("D18").Value = 396
$ws.Range("E18").Value = 5120

# Row 20: Noruega
$ws.Range("B20").Value = 4015
$ws.Range("C20").Value = 244
$ws.Range("E20").Value = 3985

# Row 25: Chequia
$ws.Range("B25").Value = 2631
$ws.Range("C25").Value = 352
$ws.Range("E25").Value = 2609

# Row 35: Rumania
$ws.Range("E35").Value = 1278
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 35

# Row 72: Bulgaria
$ws.Range("E72").Value = 313
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 7

# Row 117: Trinidad yTobago
$ws.Range("E117").Value = 70
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 3

# Row 205: San Vicente y las Granadinas
$ws.Range("D205").Value = 1
$ws.Range("E205").Value = 0

# --- Re-sorted blocks ---------------------------------------------------
# Tunez's new totals push it above Uruguay and Bosnia y Herzegovina, so the
# three rows 78-80 now hold, in order: Tunez, Uruguay, Bosnia y Herzegovina.
$ws.Range("A78").Value = "Tunez"
$ws.Range("B78").Value = 278
$ws.Range("C78").Value = 51
$ws.Range("D78").Value = 2
$ws.Range("E78").Value = 268
$ws.Range("F78").Value = 10
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 8

$ws.Range("A79").Value = "Uruguay"
$ws.Range("B79").Value = 274
$ws.Range("C79").Value = 36
$ws.Range("D79").Value = 0
$ws.Range("E79").Value = 274
$ws.Range("F79").Value = 8
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 0

$ws.Range("A80").Value = "Bosnia y Herzegovina"
$ws.Range("B80").Value = 258
$ws.Range("C80").Value = 21
$ws.Range("D80").Value = 5
$ws.Range("E80").Value = 248
$ws.Range("F80").Value = 1
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 5

# Republica de Yibuti's new totals push it above Tanzania, so rows 149-150
# now hold, in order: Republica de Yibuti, Tanzania.
$ws.Range("A149").Value = "Republica de Yibuti"
$ws.Range("B149").Value = 14
$ws.Range("C149").Value = 2
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 14
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 0

$ws.Range("A150").Value = "Tanzania"
$ws.Range("B150").Value = 14
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 13
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0
